# "added new question type: textarea"
#
# The XLSForm "survey" sheet gains a new row between the existing "text"
# question-type row and the "end group" row, introducing a "textarea"
# question type (type=textarea, name=textarea, label=Textarea). Every row
# below shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a fresh row at row 8 (just below the "text" row, just above
# "end group") — everything from row 8 down shifts to row+1.
$ws.Rows.Item(8).Insert()

# Populate the new row with the textarea question type.
$ws.Range("A8").Value = "textarea"
$ws.Range("B8").Value = "textarea"
$ws.Range("C8").Value = "Textarea"

# Match the compact row height used by the rows around it.
$ws.Rows.Item(8).RowHeight = 13.8

# Reflect the author's final cursor position.
$ws.Range("E22").Select() | Out-Null
